# Apply the grade-correction edits to the "Page 1" sheet.
# Column layout: H=Exp.20 I=PP.10 J=EX.20 K=Or.2 L=Cuest.20 M=PF.30 N=CF.105 (formula)
# N = IF(sum(H,I,J,L,M,O) < 70, IF(sum > 59, 70, sum), sum) -- recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add Cuest./PF value in M4
$ws.Range("M4").Value = 30

# Row 17: add Or. value in L17
$ws.Range("L17").Value = 14

# Row 18: add Or. value in L18
$ws.Range("L18").Value = 14

# Row 33: K33 1 -> 2, add L33
$ws.Range("K33").Value = 2
$ws.Range("L33").Value = 20

# Row 36: add L36
$ws.Range("L36").Value = 15

# Row 49: add M49
$ws.Range("M49").Value = 30

# Row 59: K59 1 -> 2, add L59
$ws.Range("K59").Value = 2
$ws.Range("L59").Value = 20

# Row 62: add M62
$ws.Range("M62").Value = 30

# Row 71: add L71
$ws.Range("L71").Value = 15

# Row 80: K80 1 -> 2, add L80
$ws.Range("K80").Value = 2
$ws.Range("L80").Value = 20

# Row 82: K82 1 -> 2, add L82
$ws.Range("K82").Value = 2
$ws.Range("L82").Value = 20

# Row 83: K83 1 -> 2, add L83
$ws.Range("K83").Value = 2
$ws.Range("L83").Value = 20

# Row 84: K84 1 -> 2, add L84
$ws.Range("K84").Value = 2
$ws.Range("L84").Value = 20

# Row 91: K91 1 -> 2, add L91
$ws.Range("K91").Value = 2
$ws.Range("L91").Value = 20

# Row 94: K94 1 -> 2, add L94
$ws.Range("K94").Value = 2
$ws.Range("L94").Value = 20

# Update the view: scroll down and move the active selection to M49
# (mirrors the author re-reviewing further down the grade sheet).
$ws.Range("M49").Select()
